$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right before "总计", based on the "2021-Q4"
#    sheet (same column layout/header "基金规模", same bold/bordered style).
# ---------------------------------------------------------------------------
$src   = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")
$src.Copy($total)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# The template only has 2 data rows (001092 / 001093); duplicate them down to
# rows 4-5 (keeping style) before overwriting with the 003719 / 161126 data.
$newSheet.Range("A2:H3").Copy($newSheet.Range("A4:H5"))
$newSheet.Range("A4").Value = 2
$newSheet.Range("A5").Value = 3

# Row 2 - 001092 (code/name already correct from the template) - new figures
$newSheet.Range("D2").Value = "'1.34"
$newSheet.Range("E2").Value = "'82.00"
$newSheet.Range("F2").Value = "'4.12"
$newSheet.Range("G2").Value = "'0.0552"
$newSheet.Range("H2").Value = 5

# Row 3 - 001093 (code/name already correct from the template) - new figures
$newSheet.Range("D3").Value = "'1.34"
$newSheet.Range("E3").Value = "'82.00"
$newSheet.Range("F3").Value = "'4.12"
$newSheet.Range("G3").Value = "'0.0552"
$newSheet.Range("H3").Value = 5

# Row 4 - 003719
$newSheet.Range("B4").Value = "'003719"
$newSheet.Range("C4").Value = "易方达标普医疗保健指数(QDII-LOF) 美元"
$newSheet.Range("D4").Value = "'0.51"
$newSheet.Range("E4").Value = "'94.20"
$newSheet.Range("F4").Value = "'1.55"
$newSheet.Range("G4").Value = "'0.0079"
$newSheet.Range("H4").Value = 6

# Row 5 - 161126
$newSheet.Range("B5").Value = "'161126"
$newSheet.Range("C5").Value = "易方达标普医疗保健指数(QDII-LOF) 人民币"
$newSheet.Range("D5").Value = "'0.51"
$newSheet.Range("E5").Value = "'94.20"
$newSheet.Range("F5").Value = "'1.55"
$newSheet.Range("G5").Value = "'0.0079"
$newSheet.Range("H5").Value = 6

# ---------------------------------------------------------------------------
# 2. Add the 2022-Q1 summary row at the top of "总计" (existing rows shift
#    down by one).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Rows(2).Insert()
$ws.Range("B2:D2").ClearFormats()

# Restore the index-column style (lost on insert) by copying from row 3.
$ws.Range("A3").Copy($ws.Range("A2"))

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q1"
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 0.13

# Renumber the shifted rows' index column (0,1,2,3,4 -> 1,2,3,4,5).
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
